$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.0498220640569395   # Accuracy
$ws1.Range("C2").Value = 0.0498220640569395   # Precision
$ws1.Range("D2").Value = 1                    # Recall
$ws1.Range("E2").Value = 0.09491525423728814  # F1 Score
$ws1.Range("F2").Value = 0.2077151335311573   # F2 Score
$ws1.Range("G2").Value = 0.5768621236133122   # F5 Score
$ws1.Range("H2").Value = 0.5430711610486891   # AUC
$ws1.Range("I2").Value = 28                   # True Positives
$ws1.Range("J2").Value = 534                  # False Positives
$ws1.Range("K2").Value = 0                    # True Negatives
$ws1.Range("L2").Value = 0                    # False Negatives

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2: class "0"
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0

# Row 3: class "1"
$ws2.Range("B3").Value = 0.0498220640569395
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.09491525423728814

# Row 4: accuracy
$ws2.Range("B4").Value = 0.0498220640569395
$ws2.Range("C4").Value = 0.0498220640569395
$ws2.Range("D4").Value = 0.0498220640569395
$ws2.Range("E4").Value = 0.0498220640569395

# Row 5: macro avg
$ws2.Range("B5").Value = 0.02491103202846975
$ws2.Range("D5").Value = 0.04745762711864407

# Row 6: weighted avg
$ws2.Range("B6").Value = 0.002482238066893783
$ws2.Range("C6").Value = 0.0498220640569395
$ws2.Range("D6").Value = 0.004728873876590867

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2: Actual 0
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 534

# Row 3: Actual 1
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
